# Re-pull / push updated "dSF" (column F) values for each data row on the
# active sheet. These values were re-derived (repull data, push all data,
# mean calculation) and differ from the previously stored (stale) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$updates = @{
    2  = 3
    3  = 1
    6  = -1
    7  = 4
    8  = -2
    9  = 1
    10 = -1
    11 = 3
    12 = 4
    13 = 1
    15 = 6
    16 = 4
    17 = 0
    18 = -1
    19 = 2
    20 = -1
    24 = -1
    25 = 4
    26 = 6
    27 = -2
    28 = -2
    29 = 1
    30 = 1
    31 = -5
    33 = -3
    34 = 3
    35 = -1
    36 = -5
    37 = -3
    38 = -2
    39 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
